$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data in row 2 and row 3 needs to be swapped for the columns that
# differ between the two records (A, B, D, E, F, G, H, Q, R). The other
# columns already hold identical values in both rows.
$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $addr2 = "$col" + "2"
    $addr3 = "$col" + "3"
    $tmp = $ws.Range($addr2).Value2
    $ws.Range($addr2).Value2 = $ws.Range($addr3).Value2
    $ws.Range($addr3).Value2 = $tmp
}
